$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.719.17"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "3.400.52"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.432"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("E9").Value = "  +6.48%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "3.396.92"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.63%  "
$ws.Range("E14").Value = "  +16.52%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "97.380.49"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "4.037.68"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +28.46%  "
$ws.Range("D19").Value = "3.405.38"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.514"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +50.72%  "
$ws.Range("E22").Value = "  +10.24%  "
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "510.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "99.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "3.579.41"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.155"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.205"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.568"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.57%  "
$ws.Range("E37").Value = "  +13.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("E39").Value = "  +11.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "524.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.857"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0424"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.81%  "
$ws.Range("E45").Value = "  -5.44%  "
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.33%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.54%  "
